# Applies the commit: the data rows (2..39) of "Avverkningsanmälningar" get
# reordered (a permutation of the existing 38 rows) and the "Förändrad"
# column (C) is bumped from 2026-02-11 (46064) to 2026-02-12 (46065) for
# every row.
#
# Strategy: snapshot every cell (A..Z) of every data row exactly as it is
# (capturing formulas verbatim and plain values with full fidelity), then
# write each snapshot back out at its new row position according to the
# permutation map below, clearing any columns the destination doesn't use.
# Finally force column C to the new date for every data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 2
$lastRow = 39
$firstCol = 1
$lastCol = 26

# new row -> old row (source of the content that should end up there)
$rowMap = @{
  2=2; 3=3; 4=4; 5=5; 6=6;
  7=10; 8=7; 9=8; 10=9;
  11=11;
  12=16; 13=15; 14=13; 15=12; 16=14;
  17=17; 18=18; 19=19;
  20=24; 21=25; 22=20; 23=21; 24=35; 25=23; 26=37; 27=30; 28=26; 29=36; 30=27;
  31=31;
  32=28; 33=34; 34=29; 35=39; 36=22; 37=32; 38=33; 39=38
}

# 1) Snapshot every cell in every data row before mutating anything.
$snapshot = @{}
for ($r = $firstRow; $r -le $lastRow; $r++) {
  $rowData = @{}
  for ($c = $firstCol; $c -le $lastCol; $c++) {
    $cell = $ws.Cells.Item($r, $c)
    if ($cell.HasFormula()) {
      $rowData[$c] = @{ "kind" = "formula"; "data" = $cell.Formula() }
    } else {
      $v = $cell.Value()
      if ($v -ne $null) {
        $rowData[$c] = @{ "kind" = "value"; "data" = $v }
      }
    }
  }
  $snapshot[$r] = $rowData
}

# 2) Write the snapshot back out in the new order, column by column,
#    clearing any cell the destination row doesn't populate.
for ($newR = $firstRow; $newR -le $lastRow; $newR++) {
  $oldR = $rowMap[$newR]
  $srcRow = $snapshot[$oldR]
  for ($c = $firstCol; $c -le $lastCol; $c++) {
    $destCell = $ws.Cells.Item($newR, $c)
    if ($srcRow.ContainsKey($c)) {
      $entry = $srcRow[$c]
      if ($entry["kind"] -eq "formula") {
        $destCell.Formula = $entry["data"]
      } else {
        $destCell.Value = $entry["data"]
      }
    } else {
      $destCell.Value = $null
    }
  }
}

# 3) Every data row's "Förändrad" (column C) moves to the new date.
for ($r = $firstRow; $r -le $lastRow; $r++) {
  $ws.Cells.Item($r, 3).Value = 46065
}
